# Applies the quant engine data refresh to the "quant Teck Fund" holdings comparison sheet.
# 1) Inserts a new "Status" column (D), shifting Jan_2026..QoQ from D:H to E:I.
# 2) Inserts two new rows (16-17) for newly exited holdings (SWIGGY LIMITED, HFCL Limited).
# 3) Refreshes every data cell with the latest values from the quant engine.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D to hold the "Status" field
$ws.Columns("D:D").Insert()

# Insert two new rows before row 16 for the new "Complete Exit" holdings
$ws.Rows("16:17").Insert()

# Header row
$ws.Range("D1").Value = "Status"
$ws.Range("E1").Value = "Jan_2026"
$ws.Range("F1").Value = "Dec_2025"
$ws.Range("G1").Value = "Oct_2025"
$ws.Range("H1").Value = "MoM"
$ws.Range("I1").Value = "QoQ"

# Row 2: Black Box Limited
$ws.Range("A2").Value = "INE676A01027"
$ws.Range("B2").Value = "Black Box Limited"
$ws.Range("C2").Value = "quant Teck Fund"
$ws.Range("D2").Value = "Adding Consistently"
$ws.Range("E2").Value = 9.991374
$ws.Range("F2").Value = 9.640183
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0.351191
$ws.Range("I2").Value = 9.991374

# Row 3: Tata Consultancy Services Limited
$ws.Range("A3").Value = "INE467B01029"
$ws.Range("B3").Value = "Tata Consultancy Services Limited"
$ws.Range("C3").Value = "quant Teck Fund"
$ws.Range("D3").Value = "Adding Consistently"
$ws.Range("E3").Value = 9.337402
$ws.Range("F3").Value = 7.453394
$ws.Range("G3").Value = 5.173243
$ws.Range("H3").Value = 1.884008000000001
$ws.Range("I3").Value = 4.164159000000001

# Row 4: Intellect Design Arena Limited
$ws.Range("A4").Value = "INE306R01017"
$ws.Range("B4").Value = "Intellect Design Arena Limited"
$ws.Range("C4").Value = "quant Teck Fund"
$ws.Range("D4").Value = "Adding Consistently"
$ws.Range("E4").Value = 7.393204
$ws.Range("F4").Value = 7.090076
$ws.Range("G4").Value = 6.017816
$ws.Range("H4").Value = 0.3031280000000001
$ws.Range("I4").Value = 1.375388

# Row 5: Digitide Solutions Limited
$ws.Range("A5").Value = "INE0U4701011"
$ws.Range("B5").Value = "Digitide Solutions Limited"
$ws.Range("C5").Value = "quant Teck Fund"
$ws.Range("D5").Value = "Reducing Consistently"
$ws.Range("E5").Value = 7.146973
$ws.Range("F5").Value = 7.42722
$ws.Range("G5").Value = 8.006432
$ws.Range("H5").Value = -0.2802470000000001
$ws.Range("I5").Value = -0.8594590000000002

# Row 6: Newgen Software Technologies Limited
$ws.Range("A6").Value = "INE619B01017"
$ws.Range("B6").Value = "Newgen Software Technologies Limited"
$ws.Range("C6").Value = "quant Teck Fund"
$ws.Range("D6").Value = "Reducing Consistently"
$ws.Range("E6").Value = 7.110316
$ws.Range("F6").Value = 9.317011
$ws.Range("G6").Value = 9.231755
$ws.Range("H6").Value = -2.206695000000001
$ws.Range("I6").Value = -2.121439

# Row 7: Adani Enterprises Limited
$ws.Range("A7").Value = "INE423A01024"
$ws.Range("B7").Value = "Adani Enterprises Limited"
$ws.Range("C7").Value = "quant Teck Fund"
$ws.Range("D7").Value = "Reducing"
$ws.Range("E7").Value = 6.212695
$ws.Range("F7").Value = 6.214914
$ws.Range("G7").Value = 6.113432
$ws.Range("H7").Value = -0.002219000000000193
$ws.Range("I7").Value = 0.09926299999999966

# Row 8: Oracle Financial Services Software Ltd
$ws.Range("A8").Value = "INE881D01027"
$ws.Range("B8").Value = "Oracle Financial Services Software Ltd"
$ws.Range("C8").Value = "quant Teck Fund"
$ws.Range("D8").Value = "Adding Consistently"
$ws.Range("E8").Value = 5.07428
$ws.Range("F8").Value = 4.535991
$ws.Range("G8").Value = 4.513687
$ws.Range("H8").Value = 0.5382889999999998
$ws.Range("I8").Value = 0.5605929999999999

# Row 9: Wipro Ltd
$ws.Range("A9").Value = "INE075A01022"
$ws.Range("B9").Value = "Wipro Ltd"
$ws.Range("C9").Value = "quant Teck Fund"
$ws.Range("D9").Value = "Reducing"
$ws.Range("E9").Value = 4.658781
$ws.Range("F9").Value = 4.672267
$ws.Range("G9").Value = 3.836748
$ws.Range("H9").Value = -0.01348599999999944
$ws.Range("I9").Value = 0.8220330000000002

# Row 10: Persistent Systems Limited
$ws.Range("A10").Value = "INE262H01021"
$ws.Range("B10").Value = "Persistent Systems Limited"
$ws.Range("C10").Value = "quant Teck Fund"
$ws.Range("D10").Value = "Fresh Entry"
$ws.Range("E10").Value = 3.990566
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 3.990566
$ws.Range("I10").Value = 3.990566

# Row 11: Sonata Software Limited
$ws.Range("A11").Value = "INE269A01021"
$ws.Range("B11").Value = "Sonata Software Limited"
$ws.Range("C11").Value = "quant Teck Fund"
$ws.Range("D11").Value = "Reducing"
$ws.Range("E11").Value = 3.136934
$ws.Range("F11").Value = 3.210402
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = -0.07346800000000009
$ws.Range("I11").Value = 3.136934

# Row 12: Infosys Limited
$ws.Range("A12").Value = "INE009A01021"
$ws.Range("B12").Value = "Infosys Limited"
$ws.Range("C12").Value = "quant Teck Fund"
$ws.Range("D12").Value = "Fresh Entry"
$ws.Range("E12").Value = 2.922662
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 2.922662
$ws.Range("I12").Value = 2.922662

# Row 13: Tata Elxsi Limited
$ws.Range("A13").Value = "INE670A01012"
$ws.Range("B13").Value = "Tata Elxsi Limited"
$ws.Range("C13").Value = "quant Teck Fund"
$ws.Range("D13").Value = "Fresh Entry"
$ws.Range("E13").Value = 2.814841
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 2.814841
$ws.Range("I13").Value = 2.814841

# Row 14: ICICI Bank Limited
$ws.Range("A14").Value = "INE090A01021"
$ws.Range("B14").Value = "ICICI Bank Limited"
$ws.Range("C14").Value = "quant Teck Fund"
$ws.Range("D14").Value = "Fresh Entry"
$ws.Range("E14").Value = 1.167371
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 1.167371
$ws.Range("I14").Value = 1.167371

# Row 15: SUN TV Network Limited
$ws.Range("A15").Value = "INE424H01027"
$ws.Range("B15").Value = "SUN TV Network Limited"
$ws.Range("C15").Value = "quant Teck Fund"
$ws.Range("D15").Value = "Complete Exit"
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 9.349346
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = -9.349346

# Row 16: SWIGGY LIMITED
$ws.Range("A16").Value = "INE00H001014"
$ws.Range("B16").Value = "SWIGGY LIMITED"
$ws.Range("C16").Value = "quant Teck Fund"
$ws.Range("D16").Value = "Complete Exit"
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 2.704044
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = -2.704044

# Row 17: HFCL Limited
$ws.Range("A17").Value = "INE548A01028"
$ws.Range("B17").Value = "HFCL Limited"
$ws.Range("C17").Value = "quant Teck Fund"
$ws.Range("D17").Value = "Complete Exit"
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 4.676026
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = -4.676026

# Row 18: R Systems International Limited
$ws.Range("A18").Value = "INE411H01032"
$ws.Range("B18").Value = "R Systems International Limited"
$ws.Range("C18").Value = "quant Teck Fund"
$ws.Range("D18").Value = "Complete Exit"
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 4.268354
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = -4.268354

# Row 19: Sasken Technologies Limited
$ws.Range("A19").Value = "INE231F01020"
$ws.Range("B19").Value = "Sasken Technologies Limited"
$ws.Range("C19").Value = "quant Teck Fund"
$ws.Range("D19").Value = "Complete Exit"
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 3.445992
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = -3.445992

# Row 20: HCL Technologies Limited
$ws.Range("A20").Value = "INE860A01027"
$ws.Range("B20").Value = "HCL Technologies Limited"
$ws.Range("C20").Value = "quant Teck Fund"
$ws.Range("D20").Value = "Complete Exit"
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 9.011688
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = -9.011688
$ws.Range("I20").Value = 0

# Row 21: Redington Limited
$ws.Range("A21").Value = "INE891D01026"
$ws.Range("B21").Value = "Redington Limited"
$ws.Range("C21").Value = "quant Teck Fund"
$ws.Range("D21").Value = "Complete Exit"
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 8.099281
$ws.Range("G21").Value = 6.77886
$ws.Range("H21").Value = -8.099281
$ws.Range("I21").Value = -6.77886
